# "process renamed to protocol in many headers"
# Rewrites the row-4 machine-readable header labels on the Enrichment
# protocol, Library preparation protocol and Sequencing protocol sheets
# so that "process_core"/"process_location" becomes
# "protocol_core"/"protocol_location" (reusing existing shared strings
# where an identical label already exists elsewhere in the workbook,
# and letting Excel create new shared strings otherwise). Also restores
# the cell selections left behind on each touched sheet plus the
# originally active sheet.

$wb = $excel.ActiveWorkbook

# --- Enrichment protocol ---------------------------------------------
$ws = $wb.Worksheets.Item("Enrichment protocol")
$ws.Activate()
$ws.Range("B4").Value = "enrichment_protocol.protocol_core.protocol_name"
$ws.Range("C4").Value = "enrichment_protocol.protocol_core.protocol_description"
$ws.Range("D4").Value = "enrichment_protocol.protocol_core.start_time"
$ws.Range("E4").Value = "enrichment_protocol.process_core.protocol_location"
$ws.Range("F4").Value = "enrichment_protocol.protocol_core.operator_identity"
$ws.Range("N4").Select()

# --- Library preparation protocol -------------------------------------
$ws = $wb.Worksheets.Item("Library preparation protocol")
$ws.Activate()
$ws.Range("A4").Value = "library_preparation_protocol.protocol_core.protocol_id"
$ws.Range("B4").Value = "library_preparation_protocol.protocol_core.protocol_name"
$ws.Range("C4").Value = "library_preparation_protocol.protocol_core.protocol_description"
$ws.Range("D4").Value = "library_preparation_protocol.protocol_core.start_time"
$ws.Range("E4").Value = "library_preparation_protocol.protocol_core.protocol_location"
$ws.Range("F4").Value = "library_preparation_protocol.protocol_core.operator_identity"
$ws.Range("AO4").Select()

# --- Sequencing protocol -----------------------------------------------
$ws = $wb.Worksheets.Item("Sequencing protocol")
$ws.Activate()
$ws.Range("A4").Value = "sequencing_protocol.protocol_core.protocol_id"
$ws.Range("B4").Value = "sequencing_protocol.protocol_core.protocol_name"
$ws.Range("C4").Value = "sequencing_protocol.protocol_core.protocol_description"
$ws.Range("D4").Value = "sequencing_protocol.protocol_core.start_time"
$ws.Range("E4").Value = "sequencing_protocol.protocol_core.protocol_location"
$ws.Range("F4").Value = "sequencing_protocol.protocol_core.operator_identity"
$ws.Range("L4").Value = "sequencing_protocol.protocol_type.text"
$ws.Range("M4").Select()

# --- Sequence files (was, and remains, the active sheet) ---------------
$ws = $wb.Worksheets.Item("Sequence files")
$ws.Activate()
$ws.Range("G4").Select()
